# Power control state diagram - update signal levels, delay and LED colours
# on the single slide of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# VSENSE_IN level labels toggle (H<->L) around the INIT/TV_ON states
$s.Shapes.Item("TextBox 49").TextFrame.TextRange.Text = "VSENSE_IN = L"
$s.Shapes.Item("TextBox 50").TextFrame.TextRange.Text = "VSENSE_IN = H"
$s.Shapes.Item("TextBox 51").TextFrame.TextRange.Text = "VSENSE_IN = L"
$s.Shapes.Item("TextBox 52").TextFrame.TextRange.Text = "VSENSE_IN = H"

# Power-off delay increased from 120s to 300s (only the first delay label)
$s.Shapes.Item("TextBox 53").TextFrame.TextRange.Text = "300s delay"

# LED indicator colours updated
$s.Shapes.Item("TextBox 59").TextFrame.TextRange.Text = "LED = BLUE"
$s.Shapes.Item("TextBox 75").TextFrame.TextRange.Text = "LED = ORANGE"

# PICONTROL_OUT levels swapped
$s.Shapes.Item("TextBox 4").TextFrame.TextRange.Text = "PICONTROL_OUT = H"
$s.Shapes.Item("TextBox 5").TextFrame.TextRange.Text = "PICONTROL_OUT = L"
